$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 5136.4
$ws.Range("I6").Value = 8416.666999999999
$ws.Range("J6").Value = 216
$ws.Range("K6").Value = 25250.001
$ws.Range("L6").Value = 648
$ws.Range("M6").Value = -25138.001
$ws.Range("N6").Value = -872
$ws.Range("H74").Value = 4101.2
$ws.Range("I74").Value = 2835.3333
$ws.Range("J74").Value = 6000
$ws.Range("K74").Value = 2835.3333
$ws.Range("L74").Value = 6000
$ws.Range("M74").Value = -1899.3333
$ws.Range("N74").Value = -7872
$ws.Range("H77").Value = 4101.2
$ws.Range("I77").Value = 2835.3333
$ws.Range("J77").Value = 6000
$ws.Range("K77").Value = 14176.6665
$ws.Range("L77").Value = 30000
$ws.Range("M77").Value = -9496.666499999999
$ws.Range("N77").Value = -39360
$ws.Range("H100").Value = 2950
$ws.Range("I100").Value = 2900
$ws.Range("K100").Value = 2900
$ws.Range("M100").Value = -2359
$ws.Range("H116").Value = 3066.0557
$ws.Range("I116").Value = 2742.3044
$ws.Range("J116").Value = 3638.8462
$ws.Range("K116").Value = 2742.3044
$ws.Range("L116").Value = 3638.8462
$ws.Range("M116").Value = 699.6956
$ws.Range("N116").Value = -10522.8462
$ws.Range("H132").Value = 16675588
$ws.Range("I132").Value = 19617564
$ws.Range("K132").Value = 58852692
$ws.Range("M132").Value = -58850162
$ws.Range("H137").Value = 1106.4517
$ws.Range("I137").Value = 979.8889
$ws.Range("J137").Value = 1281.6923
$ws.Range("K137").Value = 2939.6667
$ws.Range("L137").Value = 3845.0769
$ws.Range("M137").Value = -389.6667000000002
$ws.Range("N137").Value = -8945.0769

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2730.038
$ws.Range("I32").Value = 2461.6082
$ws.Range("J32").Value = 6702.8
$ws.Range("K32").Value = 2461.6082
$ws.Range("L32").Value = 6702.8
$ws.Range("M32").Value = -2174.6082
$ws.Range("N32").Value = -7276.8
$ws.Range("H88").Value = 2445.5454
$ws.Range("I88").Value = 1590.8
$ws.Range("J88").Value = 2696.9412
$ws.Range("K88").Value = 1590.8
$ws.Range("L88").Value = 2696.9412
$ws.Range("M88").Value = -1184.8
$ws.Range("N88").Value = -3508.9412
$ws.Range("H91").Value = 2445.5454
$ws.Range("I91").Value = 1590.8
$ws.Range("J91").Value = 2696.9412
$ws.Range("K91").Value = 1590.8
$ws.Range("L91").Value = 2696.9412
$ws.Range("M91").Value = -186.8
$ws.Range("N91").Value = -5504.9412
$ws.Range("H132").Value = 1779.1875
$ws.Range("I132").Value = 1443
$ws.Range("J132").Value = 2787.75
$ws.Range("K132").Value = 4329
$ws.Range("L132").Value = 8363.25
$ws.Range("M132").Value = -1799
$ws.Range("N132").Value = -13423.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 25074
$ws.Range("J35").Value = 25074
$ws.Range("L35").Value = 25074
$ws.Range("N35").Value = -25694
$ws.Range("H80").Value = 609.0952
$ws.Range("I80").Value = 421.66666
$ws.Range("J80").Value = 749.6667
$ws.Range("K80").Value = 421.66666
$ws.Range("L80").Value = 749.6667
$ws.Range("M80").Value = 576.33334
$ws.Range("N80").Value = -2745.6667
$ws.Range("H83").Value = 609.0952
$ws.Range("I83").Value = 421.66666
$ws.Range("J83").Value = 749.6667
$ws.Range("K83").Value = 2108.3333
$ws.Range("L83").Value = 3748.3335
$ws.Range("M83").Value = 2883.6667
$ws.Range("N83").Value = -13732.3335
$ws.Range("H113").Value = 3334933.2
$ws.Range("I113").Value = 3334933.2
$ws.Range("K113").Value = 3334933.2
$ws.Range("M113").Value = -3332763.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 560
$ws.Range("I22").Value = 383.33334
$ws.Range("J22").Value = 825
$ws.Range("K22").Value = 383.33334
$ws.Range("L22").Value = 825
$ws.Range("M22").Value = -33.33334000000002
$ws.Range("N22").Value = -1525
$ws.Range("H23").Value = 9980
$ws.Range("J23").Value = 9980
$ws.Range("L23").Value = 9980
$ws.Range("N23").Value = -10460
$ws.Range("H27").Value = 9980
$ws.Range("J27").Value = 9980
$ws.Range("L27").Value = 9980
$ws.Range("N27").Value = -10364
$ws.Range("H132").Value = 6258.391
$ws.Range("I132").Value = 8211.214
$ws.Range("J132").Value = 3220.6667
$ws.Range("K132").Value = 24633.642
$ws.Range("L132").Value = 9662.000100000001
$ws.Range("M132").Value = -22103.642
$ws.Range("N132").Value = -14722.0001
$ws.Range("H134").Value = 1889.2963
$ws.Range("I134").Value = 1951.0952
$ws.Range("J134").Value = 1673
$ws.Range("K134").Value = 5853.2856
$ws.Range("L134").Value = 5019
$ws.Range("M134").Value = -3318.2856
$ws.Range("N134").Value = -10089

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 13514696
$ws.Range("I131").Value = 166666940
$ws.Range("J131").Value = 1261.6471
$ws.Range("K131").Value = 500000820
$ws.Range("L131").Value = 3784.9413
$ws.Range("M131").Value = -499995780
$ws.Range("N131").Value = -13864.9413

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 11656.75
$ws.Range("I5").Value = 7077
$ws.Range("J5").Value = 13183.333
$ws.Range("K5").Value = 7077
$ws.Range("L5").Value = 13183.333
$ws.Range("M5").Value = -6965
$ws.Range("N5").Value = -13407.333
$ws.Range("H102").Value = 2415.8438
$ws.Range("I102").Value = 3022
$ws.Range("J102").Value = 1944.3889
$ws.Range("K102").Value = 3022
$ws.Range("L102").Value = 1944.3889
$ws.Range("M102").Value = -1400
$ws.Range("N102").Value = -5188.3889

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 213996.52
$ws.Range("I2").Value = 334201.34
$ws.Range("J2").Value = 113825.836
$ws.Range("K2").Value = 334201.34
$ws.Range("L2").Value = 113825.836
$ws.Range("M2").Value = -334089.34
$ws.Range("N2").Value = -114049.836
$ws.Range("H46").Value = 3000
$ws.Range("J46").Value = 3000
$ws.Range("L46").Value = 3000
$ws.Range("N46").Value = -3376
$ws.Range("H122").Value = 10422426
$ws.Range("I122").Value = 13895923
$ws.Range("J122").Value = 1935
$ws.Range("K122").Value = 41687769
$ws.Range("L122").Value = 5805
$ws.Range("M122").Value = -41685319
$ws.Range("N122").Value = -10705
$ws.Range("H136").Value = 5667.048
$ws.Range("I136").Value = 8007.0713
$ws.Range("J136").Value = 987
$ws.Range("K136").Value = 24021.2139
$ws.Range("L136").Value = 2961
$ws.Range("M136").Value = -21471.2139
$ws.Range("N136").Value = -8061

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 100005600
$ws.Range("J62").Value = 2500
$ws.Range("L62").Value = 2500
$ws.Range("N62").Value = -3748
$ws.Range("H65").Value = 100005600
$ws.Range("J65").Value = 2500
$ws.Range("L65").Value = 12500
$ws.Range("N65").Value = -18740
$ws.Range("H109").Value = 35377
$ws.Range("J109").Value = 35377
$ws.Range("L109").Value = 35377
$ws.Range("N109").Value = -38151
$ws.Range("H122").Value = 17335504
$ws.Range("I122").Value = 23638660
$ws.Range("J122").Value = 1824.75
$ws.Range("K122").Value = 70915980
$ws.Range("L122").Value = 5474.25
$ws.Range("M122").Value = -70913530
$ws.Range("N122").Value = -10374.25
$ws.Range("H126").Value = 90910140
$ws.Range("I126").Value = 111112040
$ws.Range("K126").Value = 333336120
$ws.Range("M126").Value = -333333650
$ws.Range("H127").Value = 63400
$ws.Range("I127").Value = 49000
$ws.Range("J127").Value = 67000
$ws.Range("K127").Value = 49000
$ws.Range("L127").Value = 67000
$ws.Range("M127").Value = -44040
$ws.Range("N127").Value = -76920
$ws.Range("H128").Value = 99990
$ws.Range("J128").Value = 99990
$ws.Range("L128").Value = 99990
$ws.Range("N128").Value = -109950
$ws.Range("H130").Value = 46250
$ws.Range("J130").Value = 46250
$ws.Range("L130").Value = 46250
$ws.Range("N130").Value = -56290
$ws.Range("H131").Value = 86660
$ws.Range("J131").Value = 86660
$ws.Range("L131").Value = 86660
$ws.Range("N131").Value = -96740
